$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.29884856315290165,
    -0.069167100505421075,
    -0.0089999996757619272,
    -0.011999999916646686,
    0.068388804065165232,
    -0.0059999996753816731,
    -0.019999999610703867,
    -0.019999999605924579,
    -0.005999999663344191,
    -0.00599999965837128,
    -0.0044999996644001783,
    -0.0059999996562836166,
    -0.0059999996508599551,
    -0.011999999622117841,
    -0.0059999996478667939,
    -0.0059999996466360006,
    -0.0059999996450423865,
    -0.0089999996311158625,
    -0.0089999996831937601,
    -0.0089999996802418991,
    -0.0089999996797631709,
    -0.0089999996794096759,
    -0.0089999996726586318,
    -0.041999999518433029,
    -0.041999999515674347,
    0.047185771513465369,
    -0.0059999996733886007,
    -0.005999999667475997,
    -0.011999999636481462,
    -0.019999999598501628,
    0.030777374747460229,
    -0.020999999590114449,
    -0.038459639069819751
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
